$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.762.16"
$ws.Range("E2").Value = "  +0.22%  "
$ws.Range("D3").Value = "2.926.20"
$ws.Range("E3").Value = "  +1.61%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "550.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.71"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +9.21%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +4.59%  "
$ws.Range("D9").Value = "2.920.78"
$ws.Range("E9").Value = "  +1.52%  "
$ws.Range("E10").Value = "  +1.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "4.76"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.78%  "
$ws.Range("E12").Value = "  +3.94%  "
$ws.Range("E13").Value = "  +3.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.122"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.70%  "
$ws.Range("D16").Value = "3.410.94"
$ws.Range("E16").Value = "  +1.67%  "
$ws.Range("E17").Value = "  +9.22%  "
$ws.Range("D18").Value = "2.923.62"
$ws.Range("E18").Value = "  +1.55%  "
$ws.Range("D19").Value = "57.740.41"
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "414.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.52%  "
$ws.Range("E21").Value = "  +4.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.693"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.91%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.14"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.63%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("E28").Value = "  +1.48%  "
$ws.Range("E29").Value = "  +6.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.32"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.90"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0978"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.64"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.934"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.88%  "
$ws.Range("D37").Value = "0.0₃0690"
$ws.Range("E37").Value = "  +12.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.25"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.74"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.60"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +12.64%  "
$ws.Range("E41").Value = "  +3.71%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0344"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.13%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "373.52"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.22%  "
$ws.Range("D44").Value = "2.694.59"
$ws.Range("E44").Value = "  +3.51%  "
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "124.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.46%  "
$ws.Range("E47").Value = "  +3.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.107"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.01%  "
$ws.Range("E49").Value = "  +1.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.98"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.87%  "
